# Apply content updates across all 7 slides of the presentation.
$p = $ppt.ActivePresentation

# Slide 1: Title slide
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(1).TextFrame.TextRange.Text = "Introduction to Graphene and 2D Nanomaterials"
$s1.Shapes.Item(2).TextFrame.TextRange.Text = "- Advances in graphene technology have spurred the synthesis of various 2D nanomaterials.`r- Transition metal oxides, metal chalcogenides, and organic compounds are key focus areas.`r- Simple and effective synthetic methods are still being pursued."

# Slide 2: Synthesis Methods for 2D Nanomaterials
$s2 = $p.Slides.Item(2)
$s2.Shapes.Item(1).TextFrame.TextRange.Text = "Synthesis Methods for 2D Nanomaterials"
$s2.Shapes.Item(2).TextFrame.TextRange.Text = "- Mechanical and liquid-phase exfoliations.`r- Ion-intercalation and exfoliation.`r- Chemical vapor deposition (CVD) and solution-phase chemical syntheses."

# Slide 3: Characterization Techniques
$s3 = $p.Slides.Item(3)
$s3.Shapes.Item(1).TextFrame.TextRange.Text = "Characterization Techniques"
$s3.Shapes.Item(2).TextFrame.TextRange.Text = "- Morphology examined using Tecnai G2 F30 S-Twin transmission electron microscope.`r- X-ray diffraction patterns recorded with Bruker D8 Advance powder X-ray diffractometer.`r- X-ray photoelectron spectra using PHI 5000 Versaprobe spectrometer."

# Slide 4: Synthesis of In4SnS8 Nanosheets
$s4 = $p.Slides.Item(4)
$s4.Shapes.Item(1).TextFrame.TextRange.Text = "Synthesis of In4SnS8 Nanosheets"
$s4.Shapes.Item(2).TextFrame.TextRange.Text = "- Synthesized via a thermal decomposition method.`r- Utilized Sn(DDTC)4 and In(DDTC)3 in OM solvent.`r- Temperature control and N2 atmosphere critical in synthesis."

# Slide 5: Adsorption and Photocatalysis Applications
$s5 = $p.Slides.Item(5)
$s5.Shapes.Item(1).TextFrame.TextRange.Text = "Adsorption and Photocatalysis Applications"
$s5.Shapes.Item(2).TextFrame.TextRange.Text = "- In4SnS8 nanosheets treated in acetic acid for surface preparation.`r- Exhibits fast adsorption and photocatalytic dual function.`r- Effective for organic dye removal in environmental remediation."

# Slide 6: Performance and Efficiency
$s6 = $p.Slides.Item(6)
$s6.Shapes.Item(1).TextFrame.TextRange.Text = "Performance and Efficiency"
$s6.Shapes.Item(2).TextFrame.TextRange.Text = "- High specific surface area of 40.34 m2 g⁻¹.`r- Superior to flower-like In4SnS8 microspheres with 24.7 m2 g⁻¹.`r- Large surface area enhances photocatalytic applications."

# Slide 7: Conclusion and Future Prospects
$s7 = $p.Slides.Item(7)
$s7.Shapes.Item(1).TextFrame.TextRange.Text = "Conclusion and Future Prospects"
$s7.Shapes.Item(2).TextFrame.TextRange.Text = "- In4SnS8 nanosheets hold potential for environmental remediation and solar energy conversion.`r- Ongoing research focuses on optimizing synthesis methods and improving efficiency.`r- Potential applications in advanced optical/electric nanodevices."
